$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Formula = "=C2+C3"
$ws.Range("C5").Select()
